# Regenerate merged AHB files:
#  - rename the "_old"/"_new" header-name suffixes to "_FV2304"/"_FV2310"
#  - turn the data range into an Excel Table (ListObject)
#  - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename header suffixes (only appear in row 1) ----------------------
$headerRange = $ws.Range("A1:U1")
$headerRange.Replace("_old", "_FV2304")
$headerRange.Replace("_new", "_FV2310")

# --- convert the used range into a native Excel table --------------------
$dataRange = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- freeze the header row (split below row 1) ---------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
